$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: D17 -> cap_Delta_17 ; definition updated to "per meg" phrasing
$ws.Range("A13").Value = "cap_Delta_17"
$ws.Range("B13").Value = "D17, triple isotopic composition of dissolved oxygen versus atmospheric O2 in per meg"

# Definitions for rows 14 and 15 updated to "per mil" phrasing (written before the
# attribute-name renames below, to match shared-string insertion order)
$ws.Range("B14").Value = "Oxygen-17 composition of dissolved oxygen versus atmospheric O2 in per mil"
$ws.Range("B15").Value = "Oxygen-18 composition of dissolved oxygen versus atmospheric O2 in per mil"

# Row 14: d17 -> delta_17
$ws.Range("A14").Value = "delta_17"

# Row 15: d18 -> delta_18
$ws.Range("A15").Value = "delta_18"

# Update the active selection to match the new state (D12:D15, active cell D12)
$ws.Range("D12:D15").Select()
